# Advent of Code 2025, day 7 — append the new runtime sample to the chart's
# backing table (row 9: Day 7, C# average runtime).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 0.00065176

# Match Excel's natural post-edit selection (the cell/range just entered).
$ws.Range("A9:B9").Select() | Out-Null
